$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "color_name" column header to "Size_name"
$ws.Range("B1").Value = "Size_name"

# Replace the "Color 1" values in column B with size values
$ws.Range("B2").Value = "1 Kg"
$ws.Range("B3").Value = "2 Kg"
$ws.Range("B4").Value = "2 Kg"
$ws.Range("B5").Value = "1 Kg"
$ws.Range("B6").Value = "2 Kg"

# Update the active selection to match the new state
$ws.Range("C4").Select()
